$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 73: convert previously-text values (stored as strings) to real
# numbers. Values stay the same, just typed as numeric now. ---
$ws.Cells.Item(73,3).Value = 0.01
$ws.Cells.Item(73,4).Value = 1
$ws.Cells.Item(73,5).Value = 0.1
$ws.Cells.Item(73,6).Value = 22
$ws.Cells.Item(73,7).Value = 24
$ws.Cells.Item(73,8).Value = 3
$ws.Cells.Item(73,9).Value = 4
$ws.Cells.Item(73,10).Value = 2409.090909090909

# --- New rows 74-80: additional median strain data, all numeric ---

# Row 74: DG011 #5, LB
$ws.Cells.Item(74,1).Value = "DG011 #5"
$ws.Cells.Item(74,2).Value = "LB"
$ws.Cells.Item(74,3).Value = 0.1
$ws.Cells.Item(74,4).Value = 0.00001
$ws.Cells.Item(74,5).Value = 0.000001
$ws.Cells.Item(74,6).Value = 177
$ws.Cells.Item(74,7).Value = 177
$ws.Cells.Item(74,8).Value = 18
$ws.Cells.Item(74,9).Value = 18
$ws.Cells.Item(74,10).Value = 177272727.2727272

# Row 75: DG012 #8, LB
$ws.Cells.Item(75,1).Value = "DG012 #8"
$ws.Cells.Item(75,2).Value = "LB"
$ws.Cells.Item(75,3).Value = 0.1
$ws.Cells.Item(75,4).Value = 0.00001
$ws.Cells.Item(75,5).Value = 0.000001
$ws.Cells.Item(75,6).Value = 255
$ws.Cells.Item(75,7).Value = 255
$ws.Cells.Item(75,8).Value = 32
$ws.Cells.Item(75,9).Value = 32
$ws.Cells.Item(75,10).Value = 260909090.9090908

# Row 76: SLM1042 #8, LB
$ws.Cells.Item(76,1).Value = "SLM1042 #8"
$ws.Cells.Item(76,2).Value = "LB"
$ws.Cells.Item(76,3).Value = 0.1
$ws.Cells.Item(76,4).Value = 0.00001
$ws.Cells.Item(76,5).Value = 0.000001
$ws.Cells.Item(76,6).Value = 232
$ws.Cells.Item(76,7).Value = 232
$ws.Cells.Item(76,8).Value = 22
$ws.Cells.Item(76,9).Value = 22
$ws.Cells.Item(76,10).Value = 230909090.9090908

# Row 77: SLM1043 #7, LB
$ws.Cells.Item(77,1).Value = "SLM1043 #7"
$ws.Cells.Item(77,2).Value = "LB"
$ws.Cells.Item(77,3).Value = 0.1
$ws.Cells.Item(77,4).Value = 0.00001
$ws.Cells.Item(77,5).Value = 0.000001
$ws.Cells.Item(77,6).Value = 139
$ws.Cells.Item(77,7).Value = 139
$ws.Cells.Item(77,8).Value = 7
$ws.Cells.Item(77,9).Value = 7
$ws.Cells.Item(77,10).Value = 132727272.7272727

# Row 78: DG011 #5, Kan
$ws.Cells.Item(78,1).Value = "DG011 #5"
$ws.Cells.Item(78,2).Value = "Kan"
$ws.Cells.Item(78,3).Value = 0.1
$ws.Cells.Item(78,4).Value = 1
$ws.Cells.Item(78,5).Value = 0.1
$ws.Cells.Item(78,6).Value = 85
$ws.Cells.Item(78,7).Value = 85
$ws.Cells.Item(78,8).Value = 11
$ws.Cells.Item(78,9).Value = 11
$ws.Cells.Item(78,10).Value = 872.7272727272726

# Row 79: DG012 #8, Kan
$ws.Cells.Item(79,1).Value = "DG012 #8"
$ws.Cells.Item(79,2).Value = "Kan"
$ws.Cells.Item(79,3).Value = 0.1
$ws.Cells.Item(79,4).Value = 1
$ws.Cells.Item(79,5).Value = 0.1
$ws.Cells.Item(79,6).Value = 175
$ws.Cells.Item(79,7).Value = 175
$ws.Cells.Item(79,8).Value = 22
$ws.Cells.Item(79,9).Value = 22
$ws.Cells.Item(79,10).Value = 1790.909090909091

# Row 80: SLM1042 #8, Kan
$ws.Cells.Item(80,1).Value = "SLM1042 #8"
$ws.Cells.Item(80,2).Value = "Kan"
$ws.Cells.Item(80,3).Value = 0.1
$ws.Cells.Item(80,4).Value = 1
$ws.Cells.Item(80,5).Value = 0.1
$ws.Cells.Item(80,6).Value = 129
$ws.Cells.Item(80,7).Value = 129
$ws.Cells.Item(80,8).Value = 8
$ws.Cells.Item(80,9).Value = 8
$ws.Cells.Item(80,10).Value = 1245.454545454545

# Row 81: SLM1043 #7, Kan -- values here stay stored as TEXT (not numbers),
# mirroring the misformatted source row that was previously at r73.
$ws.Range("C81:J81").NumberFormat = "@"
$ws.Cells.Item(81,1).Value = "SLM1043 #7"
$ws.Cells.Item(81,2).Value = "Kan"
$ws.Cells.Item(81,3).Value = "0.1"
$ws.Cells.Item(81,4).Value = "10e-1"
$ws.Cells.Item(81,5).Value = "10e-2"
$ws.Cells.Item(81,6).Value = "258"
$ws.Cells.Item(81,7).Value = "258"
$ws.Cells.Item(81,8).Value = "20"
$ws.Cells.Item(81,9).Value = "20"
$ws.Cells.Item(81,10).Value = "2527.272727272727"
